$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "roberto Santiago"
$ws.Range("B3").Value = "Afghanistan"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "+93910966393"

$ws.Range("A4").Value = "Tatiana"
$ws.Range("B4").Value = "Afghanistan"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "+93919059173"

$ws.Range("A5").Value = "roberto Santiago"
$ws.Range("B5").Value = "Afghanistan"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "+93910966392"

$ws.Range("A6").Value = "padre"
$ws.Range("B6").Value = "Afghanistan"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "+93654987245"

$ws.Range("A7").Value = "sdds"
$ws.Range("B7").Value = "Afghanistan"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "+9332323242342"

$ws.Range("A8").Value = "sdkn"
$ws.Range("B8").Value = "Andorra"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "+37632323242342"

$ws.Range("A9").Value = "teste"
$ws.Range("B9").Value = "Andorra"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "+376323232423411"

$ws.Range("A10").Value = "Thiciana Rocha "
$ws.Range("B10").Value = "Portugal"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "+351913895289"

$ws.Range("A11").Value = "2024-07-11 03:40:15"
$ws.Range("B11").Value = "roberto Santiago"
$ws.Range("C11").Value = "Afghanistan"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "+93910966393"

$ws.Range("A12").Value = "2024-07-11 03:40:15"
$ws.Range("B12").Value = "Tatiana"
$ws.Range("C12").Value = "Afghanistan"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "+93919059173"

$ws.Range("A13").Value = "2024-07-11 03:40:15"
$ws.Range("B13").Value = "roberto Santiago"
$ws.Range("C13").Value = "Afghanistan"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "+93910966392"

$ws.Range("A14").Value = "2024-07-11 03:40:15"
$ws.Range("B14").Value = "padre"
$ws.Range("C14").Value = "Afghanistan"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "+93654987245"

$ws.Range("A15").Value = "2024-07-11 03:40:15"
$ws.Range("B15").Value = "sdds"
$ws.Range("C15").Value = "Afghanistan"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "+9332323242342"

$ws.Range("A16").Value = "2024-07-11 03:40:15"
$ws.Range("B16").Value = "sdkn"
$ws.Range("C16").Value = "Andorra"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "+37632323242342"

$ws.Range("A17").Value = "2024-07-11 03:40:15"
$ws.Range("B17").Value = "teste"
$ws.Range("C17").Value = "Andorra"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "+376323232423411"

$ws.Range("A18").Value = "2024-07-11 03:40:15"
$ws.Range("B18").Value = "Thiciana Rocha "
$ws.Range("C18").Value = "Portugal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "+351913895289"
